# Updates the cryptos price/volume table with freshly scraped values
# (GitHub Actions scheduled refresh). Columns:
#   B = Coin name, C = Link, D = Price, E = Volume(1h) change
#
# Numeric-looking price strings are written with a leading apostrophe so
# Excel stores them as literal text (matching the original inlineStr cells,
# e.g. "1.00" must stay "1.00" and not collapse to the number 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.675.80"
$ws.Range("E2").Value = "  +1.79%  "

$ws.Range("D3").Value = "2.304.66"
$ws.Range("E3").Value = "  +0.77%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'316.89"
$ws.Range("E5").Value = "  +0.14%  "

$ws.Range("D6").Value = "'103.94"
$ws.Range("E6").Value = "  +0.65%  "

$ws.Range("E7").Value = "  +0.87%  "

$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("E9").Value = "  +0.63%  "

$ws.Range("D10").Value = "'39.98"
$ws.Range("E10").Value = "  +1.88%  "

$ws.Range("E11").Value = "  +0.28%  "

$ws.Range("D12").Value = "'8.52"
$ws.Range("E12").Value = "  +3.50%  "

$ws.Range("E13").Value = "  +1.59%  "

$ws.Range("D14").Value = "'0.996"
$ws.Range("E14").Value = "  +4.05%  "

$ws.Range("D15").Value = "'15.35"
$ws.Range("E15").Value = "  +1.47%  "

$ws.Range("D16").Value = "2.653.89"
$ws.Range("E16").Value = "  +0.77%  "

$ws.Range("D17").Value = "2.307.02"
$ws.Range("E17").Value = "  +0.85%  "

$ws.Range("D18").Value = "42.623.99"
$ws.Range("E18").Value = "  +1.43%  "

$ws.Range("D19").Value = "'7.64"
$ws.Range("E19").Value = "  +3.81%  "

$ws.Range("E20").Value = "  +0.75%  "

$ws.Range("D21").Value = "'13.58"
$ws.Range("E21").Value = "  +34.56%  "

$ws.Range("D22").Value = "'74.05"
$ws.Range("E22").Value = "  +1.26%  "

$ws.Range("D23").Value = "'3.54"
$ws.Range("E23").Value = "  -1.99%  "

$ws.Range("D24").Value = "'267.54"
$ws.Range("E24").Value = "  -3.51%  "

$ws.Range("E25").Value = "  -0.38%  "

$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.16%  "

$ws.Range("D27").Value = "'10.93"
$ws.Range("E27").Value = "  +1.68%  "

$ws.Range("D28").Value = "'2.35"
$ws.Range("E28").Value = "  -0.47%  "

$ws.Range("D29").Value = "'22.61"
$ws.Range("E29").Value = "  -0.59%  "

$ws.Range("D30").Value = "'38.11"
$ws.Range("E30").Value = "  +5.80%  "

$ws.Range("D31").Value = "'6.52"
$ws.Range("E31").Value = "  +12.51%  "

$ws.Range("D32").Value = "'165.93"
$ws.Range("E32").Value = "  +1.92%  "

$ws.Range("D33").Value = "'0.0883"
$ws.Range("E33").Value = "  +1.81%  "

$ws.Range("E34").Value = "  -3.25%  "

$ws.Range("D35").Value = "'2.65"
$ws.Range("E35").Value = "  -6.68%  "

$ws.Range("E36").Value = "  +0.35%  "

$ws.Range("D37").Value = "'4.60"
$ws.Range("E37").Value = "  +2.48%  "

$ws.Range("E38").Value = "  +2.38%  "

# Rows 39/40 swap ranks: NEARProtocol moves above LidoDAOToken, each with
# refreshed price/volume figures.
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D39").Value = "'3.71"
$ws.Range("E39").Value = "  -0.53%  "

$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").Value = "'2.78"
$ws.Range("E40").Value = "  -1.61%  "

$ws.Range("E41").Value = "  +14.20%  "

$ws.Range("D42").Value = "'97.63"
$ws.Range("E42").Value = "  -1.19%  "

$ws.Range("D43").Value = "'70.02"
$ws.Range("E43").Value = "  +1.46%  "

$ws.Range("E44").Value = "  +0.88%  "

$ws.Range("E45").Value = "  -0.02%  "

$ws.Range("D46").Value = "'117.53"
$ws.Range("E46").Value = "  +4.57%  "

$ws.Range("E47").Value = "  +4.03%  "

$ws.Range("D48").Value = "'80.06"
$ws.Range("E48").Value = "  +4.10%  "

$ws.Range("D49").Value = "1.649.99"
$ws.Range("E49").Value = "  +4.53%  "

$ws.Range("D50").Value = "'5.30"
$ws.Range("E50").Value = "  +0.68%  "

$ws.Range("D51").Value = "'8.89"
$ws.Range("E51").Value = "  +0.16%  "
